$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): rename / add columns ----
$ws.Range("A1").Value = "modalidade"
$ws.Range("B1").Value = "autoria_classificacao"
$ws.Range("C1").Value = "total"
$ws.Range("D1").Value = "total_sucesso"
$ws.Range("E1").Value = "particip"
$ws.Range("F1").Value = "taxa_sucesso"
$ws.Range("G1").Value = "arrecadado_sucesso"
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("P1").Value = "contribuicoes"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"
$ws.Range("U1").Value = "menor_ano"
$ws.Range("V1").Value = "maior_ano"

# New header cells (Q1:V1) need the same bold/border/centered style as the
# existing header cells -- copy format from A1 (style index 4) onto them.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("Q1:V1").PasteSpecial(-4122)

# ---- Data rows (2-6) ----
# Row 2
$ws.Range("A2").Value = "flex"
$ws.Range("B2").Value = "coletivo"
$ws.Range("C2").Value = 72
$ws.Range("D2").Value = 69
$ws.Range("E2").Value = 0.04904632152588556
$ws.Range("F2").Value = 0.9583333333333334
$ws.Range("G2").Value = 1479515.330087252
$ws.Range("H2").Value = 21442.25116068481
$ws.Range("I2").Value = 34235.40069887554
$ws.Range("J2").Value = 29.81192695893366
$ws.Range("K2").Value = 169836.9145144388
$ws.Range("L2").Value = 88.06007413874362
$ws.Range("M2").Value = 47.96955813517068
$ws.Range("N2").Value = 14.90596347946683
$ws.Range("O2").Value = 254.2443749773306
$ws.Range("P2").Value = 15501
$ws.Range("Q2").Value = 224.6521739130435
$ws.Range("R2").Value = 325.0284071787353
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 2015
$ws.Range("U2").Value = 2016
$ws.Range("V2").Value = 2023

# Row 3
$ws.Range("A3").Value = "flex"
$ws.Range("B3").Value = "empresa"
$ws.Range("C3").Value = 440
$ws.Range("D3").Value = 440
$ws.Range("E3").Value = 0.2997275204359673
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9259515.000981268
$ws.Range("H3").Value = 21044.35227495743
$ws.Range("I3").Value = 46143.04061102023
$ws.Range("J3").Value = 34.74344187043801
$ws.Range("K3").Value = 708972.7845446636
$ws.Range("L3").Value = 89.8194628153171
$ws.Range("M3").Value = 37.89772641853159
$ws.Range("N3").Value = 16.18065842403185
$ws.Range("O3").Value = 233.3973531230909
$ws.Range("P3").Value = 95943
$ws.Range("Q3").Value = 218.0522727272727
$ws.Range("R3").Value = 486.1905468340719
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 7954
$ws.Range("U3").Value = 2016
$ws.Range("V3").Value = 2023

# Row 4
$ws.Range("A4").Value = "flex"
$ws.Range("B4").Value = "feminino"
$ws.Range("C4").Value = 182
$ws.Range("D4").Value = 176
$ws.Range("E4").Value = 0.1239782016348774
$ws.Range("F4").Value = 0.967032967032967
$ws.Range("G4").Value = 1145985.994178716
$ws.Range("H4").Value = 6511.284057833613
$ws.Range("I4").Value = 6521.39877060496
$ws.Range("J4").Value = 35.53279454902379
$ws.Range("K4").Value = 29736.68915792071
$ws.Range("L4").Value = 67.58055662882595
$ws.Range("M4").Value = 23.28873679351738
$ws.Range("N4").Value = 18.47818326605706
$ws.Range("O4").Value = 154.8484188303038
$ws.Range("P4").Value = 17194
$ws.Range("Q4").Value = 97.69318181818181
$ws.Range("R4").Value = 95.76839418448722
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 453
$ws.Range("U4").Value = 2016
$ws.Range("V4").Value = 2023

# Row 5
$ws.Range("A5").Value = "flex"
$ws.Range("B5").Value = "masculino"
$ws.Range("C5").Value = 763
$ws.Range("D5").Value = 691
$ws.Range("E5").Value = 0.5197547683923706
$ws.Range("F5").Value = 0.9056356487549148
$ws.Range("G5").Value = 6465887.695217357
$ws.Range("H5").Value = 9357.290441703844
$ws.Range("I5").Value = 27421.04653337889
$ws.Range("J5").Value = 10.77163914429046
$ws.Range("K5").Value = 442290.1113560894
$ws.Range("L5").Value = 71.27543548498546
$ws.Range("M5").Value = 40.7941143515294
$ws.Range("N5").Value = 10.77163914429046
$ws.Range("O5").Value = 461.5197709071476
$ws.Range("P5").Value = 74806
$ws.Range("Q5").Value = 108.2575976845152
$ws.Range("R5").Value = 214.0671103886876
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 3474
$ws.Range("U5").Value = 2016
$ws.Range("V5").Value = 2023

# Row 6
$ws.Range("A6").Value = "flex"
$ws.Range("B6").Value = "outros"
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 0.007493188010899182
$ws.Range("F6").Value = 0.6363636363636364
$ws.Range("G6").Value = 11227.91709450537
$ws.Range("H6").Value = 1603.98815635791
$ws.Range("I6").Value = 2112.496775736096
$ws.Range("J6").Value = 42.35779660756832
$ws.Range("K6").Value = 5515.844600589859
$ws.Range("L6").Value = 45.24320624776205
$ws.Range("M6").Value = 14.93011251059404
$ws.Range("N6").Value = 21.17889830378416
$ws.Range("O6").Value = 63.40051265045815
$ws.Range("P6").Value = 202
$ws.Range("Q6").Value = 28.85714285714286
$ws.Range("R6").Value = 34.72956747704838
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 87
$ws.Range("U6").Value = 2017
$ws.Range("V6").Value = 2023

# ---- Number formats per column (matches style classes 1/2/3; General = default) ----
$ws.Range("C2:C6").NumberFormat = "#,##0"
$ws.Range("D2:D6").NumberFormat = "#,##0"
$ws.Range("E2:E6").NumberFormat = "0.00%"
$ws.Range("F2:F6").NumberFormat = "0.00%"
$ws.Range("G2:G6").NumberFormat = "R$ #,##0.00"
$ws.Range("H2:H6").NumberFormat = "R$ #,##0.00"
$ws.Range("I2:I6").NumberFormat = "R$ #,##0.00"
$ws.Range("J2:J6").NumberFormat = "R$ #,##0.00"
$ws.Range("K2:K6").NumberFormat = "R$ #,##0.00"
$ws.Range("L2:L6").NumberFormat = "R$ #,##0.00"
$ws.Range("M2:M6").NumberFormat = "R$ #,##0.00"
$ws.Range("N2:N6").NumberFormat = "R$ #,##0.00"
$ws.Range("O2:O6").NumberFormat = "R$ #,##0.00"
$ws.Range("P2:P6").NumberFormat = "#,##0"
$ws.Range("Q2:Q6").NumberFormat = "#,##0"
$ws.Range("R2:R6").NumberFormat = "#,##0"
$ws.Range("S2:S6").NumberFormat = "#,##0"
$ws.Range("T2:T6").NumberFormat = "#,##0"
